$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on the worksheet's used range
$lastRow = $ws.UsedRange.Rows.Count

# Column C holds the "Förändrad" (changed) date. Every data row (2..lastRow)
# currently stores serial date 45171 (2023-09-02) and should be bumped by
# one day to 45172 (2023-09-03).
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }
}
